# Feature - new update
# Populate the "email" sheet with the harvested OAuth credential record:
# a header row (A1:I1) plus a single data row (A2:I2). All values are
# plain text except the expiry_date (I2), which is a numeric epoch-ms
# timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "email"
$ws.Cells.Item(1, 2).Value = "client_id"
$ws.Cells.Item(1, 3).Value = "client_secret"
$ws.Cells.Item(1, 4).Value = "redirect_uris"
$ws.Cells.Item(1, 5).Value = "access_token"
$ws.Cells.Item(1, 6).Value = "refresh_token"
$ws.Cells.Item(1, 7).Value = "scope"
$ws.Cells.Item(1, 8).Value = "token_type"
$ws.Cells.Item(1, 9).Value = "expiry_date"

# Data row
$ws.Cells.Item(2, 1).Value = "alyssamarie69554@gmail.com"
$ws.Cells.Item(2, 2).Value = "440677329044-5tb7ovk4h5gs19ktcfau55hqf192o86q.apps.googleusercontent.com"
$ws.Cells.Item(2, 3).Value = "I6UAzzKNvJpBIG4MXAwpGh9U"
$ws.Cells.Item(2, 4).Value = "urn:ietf:wg:oauth:2.0:oob"
$ws.Cells.Item(2, 5).Value = "ya29.Il-EB1bTQsiCmYY8ve-xKDvVtXAmXLvZ9kt9InDE1wnk_7TUrqOAznuISim3c2iCfP7g9WXBFOrhdMXVS7JV3O89OShIJhvA82eYZ7vdut0ge3ZI1EUwPZl--D37px1-sg"
$ws.Cells.Item(2, 6).Value = "1/PCbHoOKHgrk7f7moy_GFTUeKg8ZHYw6deBERYKMkhyg"
$ws.Cells.Item(2, 7).Value = "https://mail.google.com/"
$ws.Cells.Item(2, 8).Value = "Bearer"
$ws.Cells.Item(2, 9).Value = 1568523661006
